$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the existing "District" column (F), shifting it to G.
$ws.Columns.Item(6).Insert()

# New header for the inserted "Address" column.
$ws.Range("F2").Value = "Address"

# Populate the "Address" column per row, derived from the school/taluk portion
# of the NAMES (B) column text (the part between the teacher's name and the
# trailing district, which already lived in the old "District" column).
    $ws.Range("F3").Value = "Morarji Desai Residential School MadanabaviHonnali"
    $ws.Range("F4").Value = "G H S KanchugaranahalliChannagiri"
    $ws.Range("F5").Value = "Govt. High School GoverahalliHarapanahalli"
    $ws.Range("F6").Value = "G H S ParagoduBagepalli"
    $ws.Range("F7").Value = "Vivekananda High School NyamathiHonnali"
    $ws.Range("F8").Value = "G G H S Channagiri"
    $ws.Range("F9").Value = "S R R H SchoolH M HalliJagalur"
    $ws.Range("F10").Value = "S V B H S KyasinakereHonnali"
    $ws.Range("F11").Value = "S N H S ChigateriHarapanahalli"
    $ws.Range("F12").Value = "G G H S Bagepalli"
    $ws.Range("F13").Value = "Sri Hucheshwar High School KamatgiHungund"
    $ws.Range("F14").Value = "Govt. P U College ChikkalagundiBilagi"
    $ws.Range("F15").Value = "Govt. High School AnawalBadami"
    $ws.Range("F16").Value = "Hema vema High SchoolUgalawatBadami"
    $ws.Range("F17").Value = "Govt. High School BevinamattiHungunda"
    $ws.Range("F18").Value = "G H P S Chinchalakatti"
    $ws.Range("F20").Value = "S G H S TumbaHungund"
    $ws.Range("F21").Value = "H P S BinjawadagiHungund"
    $ws.Range("F22").Value = "G H S ThimmampalliBagepalli"
    $ws.Range("F23").Value = "G H S DarinayakanapalyaGowribidanur"
    $ws.Range("F24").Value = "G H S TolamattiBilgi"
    $ws.Range("F25").Value = "Govt. High School BillurBagepalli"
    $ws.Range("F26").Value = "G H S BadagandiBilagi"
    $ws.Range("F27").Value = "Govt. High School H GopagondanahalliHonnali"
    $ws.Range("F28").Value = "Govt. High School ChakaveluBagepalli"
    $ws.Range("F29").Value = "National High School BelaguttiHonnali"
    $ws.Range("F30").Value = "S A H S KondadahalliChannagiri"
    $ws.Range("F31").Value = "Govt. High SchoolHirebadawadagiHungund"
    $ws.Range("F32").Value = "S S P U CollegeBilgi"
    $ws.Range("F33").Value = "Sri Maruthi PragathiHigh School MudaloduGowribidanur"
    $ws.Range("F34").Value = "G H S MuragamaleChintamani"
    $ws.Range("F35").Value = "K N C S S U School"
    $ws.Range("F36").Value = "Harapanahalli"
    $ws.Range("F37").Value = "G Hanumath Reddy Memorial Residential High SchoolJagalur town Jagalur"
    $ws.Range("F38").Value = "S S R High SchoolBannikudoHonnali"
    $ws.Range("F39").Value = "G H S ChiradoniChannagiri"
    $ws.Range("F40").Value = "S R M P P G H S MuthigiHarapanahalli"
    $ws.Range("F41").Value = "Sree Bheemeshwara High School YidagurGauribidanur"
    $ws.Range("F42").Value = "S TJ G J C HirekogalurChannagiri"
    $ws.Range("F43").Value = "Shree S B H S SulebhaviHunagund"
    $ws.Range("F44").Value = "S V P H S KatapurBadami"
    $ws.Range("F45").Value = "G J C KariganurKathalagereChannagiri"
    $ws.Range("F46").Value = "G H S YaralakatteJagalur"
    $ws.Range("F47").Value = "Adarsha Vidyalaya (R M S A) Bagepalli"
    $ws.Range("F48").Value = "Pathi Adinaranaiah Ramalakshmamma High School MunganahalliChintamani"
    $ws.Range("F49").Value = "T K G M H S T GollahalliChintamani"
    $ws.Range("F50").Value = "G H S PathapalyaBagepalli"
    $ws.Range("F52").Value = "Chintamani"
    $ws.Range("F53").Value = "Upgraded Govt. High SchoolHoovinhalliHungud"
    $ws.Range("F54").Value = "Venktesh High School GuledaguddaBadami"
    $ws.Range("F55").Value = "G H S AmalazariBilagi"
    $ws.Range("F56").Value = "G H S Kolur (RC)Bilagi"
    $ws.Range("F57").Value = "Shri Kanchaneshwari H S GuledguddaBadami"
    $ws.Range("F58").Value = "G H S NamagundluGowribidanur"
    $ws.Range("F59").Value = "G H S YagavakoteChintamani"
    $ws.Range("F60").Value = "Sree Gurusiddeshwara High School BhidaraqereJagalur"
    $ws.Range("F61").Value = "H P S ChinnapurHungund"
    $ws.Range("F62").Value = "G H S KainakattiBadami"
    $ws.Range("F63").Value = "G H S SonaganahalliGauribidanur"
    $ws.Range("F64").Value = "G H P School NaraspurBadami"
    $ws.Range("F65").Value = "G H S ManchenahalliGauribidanur"
    $ws.Range("F66").Value = "Sri Ramakrishna Vidya Vardhaka High SchoolYenigadaleChintamani"
    $ws.Range("F67").Value = "M L A M P S No 1 Badami"
    $ws.Range("F70").Value = "Govt. Girls High School Gudibande"
    $ws.Range("F71").Value = "G H S MinakanagurkiGauribidanur"
    $ws.Range("F72").Value = "Govt. High School JambaladinniHungund"
    $ws.Range("F73").Value = "Govt. High SchoolMittemariBagepalli"
    $ws.Range("F74").Value = "G H S HarakanaluHarapanahalli"
    $ws.Range("F75").Value = "Govt. P U College (H S)IlkalHunagund"

